$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows to append (duplicating existing rows 5,8,4,6,7,3,2, in that
# order) as rows 9-15, matching the scraped data refresh described by the
# commit.
$batsman = "Sanju Samson" + [char]0x00A0

$newRows = @(
    @(" Dubai (DSC)",  " October 14 2020",   "Capitals won by 13 runs",                           "Rajasthan Royals", "Delhi Capitals",              $batsman, "25", "18", "0", "2", "138.88"),
    @(" Dubai (DSC)",  " September 30 2020", "KKR won by 37 runs",                                "Rajasthan Royals", "Kolkata Knight Riders",       $batsman, "8",  "9",  "1", "0", "88.88"),
    @(" Sharjah",      " September 27 2020", "Royals won by 4 wickets (with 3 balls remaining)",  "Rajasthan Royals", "Kings XI Punjab",             $batsman, "85", "42", "4", "7", "202.38"),
    @(" Dubai (DSC)",  " October 11 2020",   "Royals won by 5 wickets (with 1 ball remaining)",   "Rajasthan Royals", "Sunrisers Hyderabad",         $batsman, "26", "25", "3", "0", "104.00"),
    @(" Abu Dhabi",    " October 06 2020",   "Mumbai won by 57 runs",                             "Rajasthan Royals", "Mumbai Indians",              $batsman, "0",  "3",  "0", "0", "0.00"),
    @(" Abu Dhabi",    " October 03 2020",   "RCB won by 8 wickets (with 5 balls remaining)",     "Rajasthan Royals", "Royal Challengers Bangalore", $batsman, "4",  "3",  "1", "0", "133.33"),
    @(" Sharjah",      " October 09 2020",   "Capitals won by 46 runs",                           "Rajasthan Royals", "Delhi Capitals",              $batsman, "5",  "9",  "0", "0", "55.55")
)

$startRow = 9
$endRow = $startRow + $newRows.Count - 1

# The source data (scraped JSON) stores every column as text, including the
# numeric-looking ones (runs/balls/4s/6s/strike-rate) - force text storage
# so "25" stays "25" and not the number 25.
$ws.Range("A$startRow`:K$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Output "Appended $($newRows.Count) rows starting at row $startRow"
